$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(8, 8).Value = 3853.889
$ws.Cells.Item(8, 9).Value = 2054.4443
$ws.Cells.Item(8, 10).Value = 5653.3335
$ws.Cells.Item(8, 11).Value = 6163.3329
$ws.Cells.Item(8, 12).Value = 16960.0005
$ws.Cells.Item(8, 13).Value = -6024.3329
$ws.Cells.Item(8, 14).Value = -17238.0005

$ws.Cells.Item(39, 8).Value = 400.5
$ws.Cells.Item(39, 9).Value = 175.75
$ws.Cells.Item(39, 10).Value = 512.875
$ws.Cells.Item(39, 11).Value = 527.25
$ws.Cells.Item(39, 12).Value = 1538.625
$ws.Cells.Item(39, 13).Value = -231.25
$ws.Cells.Item(39, 14).Value = -2130.625

$ws.Cells.Item(41, 8).Value = 1093.8572
$ws.Cells.Item(41, 9).Value = 719.9091
$ws.Cells.Item(41, 11).Value = 719.9091
$ws.Cells.Item(41, 13).Value = -279.9091

$ws.Cells.Item(74, 8).Value = 28236.75
$ws.Cells.Item(74, 9).Value = 28236.75
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 28236.75
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = ""
$ws.Cells.Item(74, 14).Value = -27300.75

$ws.Cells.Item(77, 8).Value = 28236.75
$ws.Cells.Item(77, 9).Value = 28236.75
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 141183.75
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = ""
$ws.Cells.Item(77, 14).Value = -136503.75

$ws.Cells.Item(113, 8).Value = 38467384
$ws.Cells.Item(113, 9).Value = 10006002
$ws.Cells.Item(113, 10).Value = 133338664
$ws.Cells.Item(113, 11).Value = 10006002
$ws.Cells.Item(113, 12).Value = 133338664
$ws.Cells.Item(113, 13).Value = -10002748
$ws.Cells.Item(113, 14).Value = -133345172

$ws.Cells.Item(135, 8).Value = 11653.3125
$ws.Cells.Item(135, 9).Value = 3998
$ws.Cells.Item(135, 11).Value = 35982
$ws.Cells.Item(135, 13).Value = -33447

$ws.Cells.Item(138, 8).Value = 1917.225
$ws.Cells.Item(138, 10).Value = 3416.389
$ws.Cells.Item(138, 12).Value = 10249.167
$ws.Cells.Item(138, 14).Value = -20529.167

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(7, 8).Value = 68000
$ws.Cells.Item(7, 10).Value = 68000
$ws.Cells.Item(7, 12).Value = 68000
$ws.Cells.Item(7, 14).Value = -68228

$ws.Cells.Item(74, 8).Value = 31271810
$ws.Cells.Item(74, 9).Value = 62502504
$ws.Cells.Item(74, 11).Value = 62502504
$ws.Cells.Item(74, 13).Value = -62501630

$ws.Cells.Item(77, 8).Value = 31271810
$ws.Cells.Item(77, 9).Value = 62502504
$ws.Cells.Item(77, 11).Value = 312512520
$ws.Cells.Item(77, 13).Value = -312508152

$ws.Cells.Item(132, 8).Value = 6686.6665
$ws.Cells.Item(132, 9).Value = 3542.2942
$ws.Cells.Item(132, 11).Value = 10626.8826
$ws.Cells.Item(132, 13).Value = -8096.882599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(75, 8).Value = 15000
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 15000
$ws.Cells.Item(75, 11).Value = 0
$ws.Cells.Item(75, 12).Value = ""
$ws.Cells.Item(75, 13).Value = 15000
$ws.Cells.Item(75, 14).Value = -16872

$ws.Cells.Item(78, 8).Value = 15000
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 15000
$ws.Cells.Item(78, 11).Value = 0
$ws.Cells.Item(78, 12).Value = ""
$ws.Cells.Item(78, 13).Value = 45000
$ws.Cells.Item(78, 14).Value = -54360

$ws.Cells.Item(82, 8).Value = 29666.5
$ws.Cells.Item(82, 9).Value = 9333
$ws.Cells.Item(82, 11).Value = 9333
$ws.Cells.Item(82, 13).Value = -8950

$ws.Cells.Item(85, 8).Value = 29666.5
$ws.Cells.Item(85, 9).Value = 9333
$ws.Cells.Item(85, 11).Value = 9333
$ws.Cells.Item(85, 13).Value = -8007

$ws.Cells.Item(107, 8).Value = 1693.75
$ws.Cells.Item(107, 9).Value = 1463.1177
$ws.Cells.Item(107, 11).Value = 1463.1177
$ws.Cells.Item(107, 13).Value = 456.8823

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(29, 8).Value = 100
$ws.Cells.Item(29, 9).Value = 100
$ws.Cells.Item(29, 11).Value = 100
$ws.Cells.Item(29, 13).Value = 193

$ws.Cells.Item(31, 8).Value = 810795
$ws.Cells.Item(31, 9).Value = 11286.934
$ws.Cells.Item(31, 10).Value = 1667410.8
$ws.Cells.Item(31, 11).Value = 11286.934
$ws.Cells.Item(31, 12).Value = 1667410.8
$ws.Cells.Item(31, 13).Value = -10991.934
$ws.Cells.Item(31, 14).Value = -1668000.8

$ws.Cells.Item(34, 8).Value = 810795
$ws.Cells.Item(34, 9).Value = 11286.934
$ws.Cells.Item(34, 10).Value = 1667410.8
$ws.Cells.Item(34, 11).Value = 11286.934
$ws.Cells.Item(34, 12).Value = 1667410.8
$ws.Cells.Item(34, 13).Value = -11084.934
$ws.Cells.Item(34, 14).Value = -1667814.8

$ws.Cells.Item(107, 8).Value = 945.05
$ws.Cells.Item(107, 9).Value = 661.3333
$ws.Cells.Item(107, 11).Value = 661.3333
$ws.Cells.Item(107, 13).Value = 1258.6667

$ws.Cells.Item(118, 8).Value = 87979
$ws.Cells.Item(118, 10).Value = 87979
$ws.Cells.Item(118, 12).Value = 87979
$ws.Cells.Item(118, 14).Value = -91293

$ws.Cells.Item(132, 8).Value = 2241.2222
$ws.Cells.Item(132, 9).Value = 2258.72
$ws.Cells.Item(132, 11).Value = 6776.16
$ws.Cells.Item(132, 13).Value = -4246.16

$ws.Cells.Item(134, 8).Value = 458035.88
$ws.Cells.Item(134, 9).Value = 626549.5600000001
$ws.Cells.Item(134, 11).Value = 1879648.68
$ws.Cells.Item(134, 13).Value = -1877113.68

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(10, 8).Value = 16.90909
$ws.Cells.Item(10, 9).Value = 17.6
$ws.Cells.Item(10, 11).Value = 52.8
$ws.Cells.Item(10, 13).Value = 86.19999999999999

$ws.Cells.Item(12, 8).Value = 366087.47
$ws.Cells.Item(12, 10).Value = 594188.9
$ws.Cells.Item(12, 12).Value = 1782566.7
$ws.Cells.Item(12, 14).Value = -1782912.7

$ws.Cells.Item(50, 8).Value = 456.81818
$ws.Cells.Item(50, 10).Value = 461.90475
$ws.Cells.Item(50, 12).Value = 1385.71425
$ws.Cells.Item(50, 14).Value = -2347.71425

$ws.Cells.Item(53, 8).Value = 456.81818
$ws.Cells.Item(53, 10).Value = 461.90475
$ws.Cells.Item(53, 12).Value = 1385.71425
$ws.Cells.Item(53, 14).Value = -2347.71425

$ws.Cells.Item(125, 8).Value = 15005.5
$ws.Cells.Item(125, 10).Value = 15005.5
$ws.Cells.Item(125, 12).Value = 45016.5
$ws.Cells.Item(125, 14).Value = -54856.5

$ws.Cells.Item(126, 8).Value = 5111.1113
$ws.Cells.Item(126, 10).Value = 10000
$ws.Cells.Item(126, 12).Value = 30000
$ws.Cells.Item(126, 14).Value = -39880

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 18332.9
$ws.Cells.Item(80, 9).Value = 14666.134
$ws.Cells.Item(80, 11).Value = 14666.134
$ws.Cells.Item(80, 13).Value = -13668.134

$ws.Cells.Item(83, 8).Value = 18332.9
$ws.Cells.Item(83, 9).Value = 14666.134
$ws.Cells.Item(83, 11).Value = 73330.67
$ws.Cells.Item(83, 13).Value = -68338.67

$ws.Cells.Item(107, 8).Value = 1473.25
$ws.Cells.Item(107, 9).Value = 1272.125
$ws.Cells.Item(107, 11).Value = 1272.125
$ws.Cells.Item(107, 13).Value = 647.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(10, 8).Value = 6245.75
$ws.Cells.Item(10, 9).Value = 4997.5
$ws.Cells.Item(10, 10).Value = 7494
$ws.Cells.Item(10, 11).Value = 4997.5
$ws.Cells.Item(10, 12).Value = 7494
$ws.Cells.Item(10, 13).Value = -4857.5
$ws.Cells.Item(10, 14).Value = -7774

$ws.Cells.Item(45, 8).Value = 5041
$ws.Cells.Item(45, 9).Value = 5041
$ws.Cells.Item(45, 11).Value = 5041
$ws.Cells.Item(45, 13).Value = -4634

$ws.Cells.Item(46, 8).Value = 2963.923
$ws.Cells.Item(46, 10).Value = 4221.3335
$ws.Cells.Item(46, 12).Value = 4221.3335
$ws.Cells.Item(46, 14).Value = -4597.3335

$ws.Cells.Item(60, 8).Value = 50000
$ws.Cells.Item(60, 10).Value = 50000
$ws.Cells.Item(60, 12).Value = 50000
$ws.Cells.Item(60, 14).Value = -51018

$ws.Cells.Item(122, 8).Value = 5943.8184
$ws.Cells.Item(122, 9).Value = 5264.6665
$ws.Cells.Item(122, 11).Value = 15793.9995
$ws.Cells.Item(122, 13).Value = -13343.9995

$ws.Cells.Item(132, 8).Value = 57730.082
$ws.Cells.Item(132, 9).Value = 31567
$ws.Cells.Item(132, 11).Value = 94701
$ws.Cells.Item(132, 13).Value = -92171

$ws.Cells.Item(136, 8).Value = 102683.48
$ws.Cells.Item(136, 9).Value = 60755.766
$ws.Cells.Item(136, 11).Value = 182267.298
$ws.Cells.Item(136, 13).Value = -179717.298

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 11).Value = 0
$ws.Cells.Item(81, 13).Value = ""

$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 11).Value = 0
$ws.Cells.Item(84, 13).Value = ""

$ws.Cells.Item(107, 8).Value = 33334684
$ws.Cells.Item(107, 9).Value = 45456104
$ws.Cells.Item(107, 11).Value = 136368312
$ws.Cells.Item(107, 13).Value = -136366392

$ws.Cells.Item(125, 8).Value = 88997.25
$ws.Cells.Item(125, 10).Value = 88997.25
$ws.Cells.Item(125, 12).Value = 88997.25
$ws.Cells.Item(125, 14).Value = -98837.25

$ws.Cells.Item(132, 8).Value = 2688.3125
$ws.Cells.Item(132, 9).Value = 1975.9
$ws.Cells.Item(132, 11).Value = 5927.700000000001
$ws.Cells.Item(132, 13).Value = -3397.700000000001
